# Applies the scheduled-runner data refresh described in the commit.
# Each Leve row's market-price / profit columns (H.. N) are updated in
# place to the freshly pulled Universalis price data; a couple of rows
# (GSM!75 and GSM!78) lost their HQ-profit figure entirely because the
# refreshed HQ price came back as zero, so LeveProfitHQ (N) is cleared.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9805790
$ws.Range("I137").Value = 15153474
$ws.Range("J137").Value = 1702.8889
$ws.Range("K137").Value = 45460422
$ws.Range("L137").Value = 5108.6667
$ws.Range("M137").Value = -45457872
$ws.Range("N137").Value = -10208.6667

$ws.Range("H138").Value = 1608.9818
$ws.Range("J138").Value = 1937.8572
$ws.Range("L138").Value = 5813.571599999999
$ws.Range("N138").Value = -16093.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5966765
$ws.Range("I32").Value = 8784.093000000001
$ws.Range("J32").Value = 25673934
$ws.Range("K32").Value = 8784.093000000001
$ws.Range("L32").Value = 25673934
$ws.Range("M32").Value = -8497.093000000001
$ws.Range("N32").Value = -25674508

$ws.Range("H61").Value = 2060.9111
$ws.Range("I61").Value = 1174.4849
$ws.Range("J61").Value = 4498.5835
$ws.Range("K61").Value = 1174.4849
$ws.Range("L61").Value = 4498.5835
$ws.Range("M61").Value = -962.4848999999999
$ws.Range("N61").Value = -4922.5835

$ws.Range("H63").Value = 3937.8235
$ws.Range("I63").Value = 1924.6428
$ws.Range("J63").Value = 13332.667
$ws.Range("K63").Value = 1924.6428
$ws.Range("L63").Value = 13332.667
$ws.Range("M63").Value = -1238.6428
$ws.Range("N63").Value = -14704.667

$ws.Range("H66").Value = 3937.8235
$ws.Range("I66").Value = 1924.6428
$ws.Range("J66").Value = 13332.667
$ws.Range("K66").Value = 9623.214
$ws.Range("L66").Value = 66663.33499999999
$ws.Range("M66").Value = -6191.214
$ws.Range("N66").Value = -73527.33499999999

$ws.Range("H97").Value = 295.23077
$ws.Range("I97").Value = 223.8
$ws.Range("J97").Value = 533.3333
$ws.Range("K97").Value = 223.8
$ws.Range("L97").Value = 533.3333
$ws.Range("M97").Value = 272.2
$ws.Range("N97").Value = -1525.3333

$ws.Range("H132").Value = 1841370.9
$ws.Range("I132").Value = 1283.2128
$ws.Range("J132").Value = 5959662
$ws.Range("K132").Value = 3849.6384
$ws.Range("L132").Value = 17878986
$ws.Range("M132").Value = -1319.6384
$ws.Range("N132").Value = -17884046

$ws.Range("H136").Value = 2060.9111
$ws.Range("I136").Value = 1174.4849
$ws.Range("J136").Value = 4498.5835
$ws.Range("K136").Value = 3523.4547
$ws.Range("L136").Value = 13495.7505
$ws.Range("M136").Value = -973.4546999999998
$ws.Range("N136").Value = -18595.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 749.9167
$ws.Range("I94").Value = 749.9
$ws.Range("J94").Value = 750
$ws.Range("K94").Value = 749.9
$ws.Range("L94").Value = 750
$ws.Range("M94").Value = -298.9
$ws.Range("N94").Value = -1652

$ws.Range("H112").Value = 36489.668
$ws.Range("J112").Value = 36489.668
$ws.Range("L112").Value = 36489.668
$ws.Range("N112").Value = -39443.668

$ws.Range("H134").Value = 4137.5576
$ws.Range("I134").Value = 1748.1892
$ws.Range("J134").Value = 7821.1665
$ws.Range("K134").Value = 5244.5676
$ws.Range("L134").Value = 23463.4995
$ws.Range("M134").Value = -2709.5676
$ws.Range("N134").Value = -28533.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2126.7727
$ws.Range("I16").Value = 2339.9333
$ws.Range("J16").Value = 1670
$ws.Range("K16").Value = 2339.9333
$ws.Range("L16").Value = 1670
$ws.Range("M16").Value = -2052.9333
$ws.Range("N16").Value = -2244

$ws.Range("H31").Value = 7144890.5
$ws.Range("I31").Value = 1330.7858
$ws.Range("J31").Value = 17860230
$ws.Range("K31").Value = 1330.7858
$ws.Range("L31").Value = 17860230
$ws.Range("M31").Value = -1035.7858
$ws.Range("N31").Value = -17860820

$ws.Range("H34").Value = 7144890.5
$ws.Range("I34").Value = 1330.7858
$ws.Range("J34").Value = 17860230
$ws.Range("K34").Value = 1330.7858
$ws.Range("L34").Value = 17860230
$ws.Range("M34").Value = -1128.7858
$ws.Range("N34").Value = -17860634

$ws.Range("H105").Value = 477282.72
$ws.Range("I105").Value = 667709.3
$ws.Range("J105").Value = 1216.1666
$ws.Range("K105").Value = 667709.3
$ws.Range("L105").Value = 1216.1666
$ws.Range("M105").Value = -665962.3
$ws.Range("N105").Value = -4710.1666

$ws.Range("H113").Value = 2126.7727
$ws.Range("I113").Value = 2339.9333
$ws.Range("J113").Value = 1670
$ws.Range("K113").Value = 2339.9333
$ws.Range("L113").Value = 1670
$ws.Range("M113").Value = -169.9333000000001
$ws.Range("N113").Value = -6010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H97").Value = 1486.1111
$ws.Range("I97").Value = 1546.875
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1546.875
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -1050.875
$ws.Range("N97").Value = -1992

$ws.Range("H132").Value = 993748.2
$ws.Range("I132").Value = 1264046.2
$ws.Range("J132").Value = 2655.3333
$ws.Range("K132").Value = 3792138.6
$ws.Range("L132").Value = 7965.999899999999
$ws.Range("M132").Value = -3789608.6
$ws.Range("N132").Value = -13025.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2701.6296
$ws.Range("I68").Value = 2468.7334
$ws.Range("J68").Value = 2992.75
$ws.Range("K68").Value = 2468.7334
$ws.Range("L68").Value = 2992.75
$ws.Range("M68").Value = -1719.7334
$ws.Range("N68").Value = -4490.75

$ws.Range("H71").Value = 2701.6296
$ws.Range("I71").Value = 2468.7334
$ws.Range("J71").Value = 2992.75
$ws.Range("K71").Value = 12343.667
$ws.Range("L71").Value = 14963.75
$ws.Range("M71").Value = -8599.667000000001
$ws.Range("N71").Value = -22451.75

$ws.Range("H93").Value = 239544.05
$ws.Range("I93").Value = 313773.2
$ws.Range("J93").Value = 2010.8
$ws.Range("K93").Value = 313773.2
$ws.Range("L93").Value = 2010.8
$ws.Range("M93").Value = -312525.2
$ws.Range("N93").Value = -4506.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4386
$ws.Range("I62").Value = 3217.625
$ws.Range("J62").Value = 6722.75
$ws.Range("K62").Value = 3217.625
$ws.Range("L62").Value = 6722.75
$ws.Range("M62").Value = -2593.625
$ws.Range("N62").Value = -7970.75

$ws.Range("H65").Value = 4386
$ws.Range("I65").Value = 3217.625
$ws.Range("J65").Value = 6722.75
$ws.Range("K65").Value = 16088.125
$ws.Range("L65").Value = 33613.75
$ws.Range("M65").Value = -12968.125
$ws.Range("N65").Value = -39853.75

$ws.Range("H136").Value = 8073960
$ws.Range("I136").Value = 9625934
$ws.Range("J136").Value = 3695
$ws.Range("K136").Value = 28877802
$ws.Range("L136").Value = 11085
$ws.Range("M136").Value = -28875252
$ws.Range("N136").Value = -16185
